# Deploy the implementation guide.
# - Rename the "Include from Ferlab.bio CodeS" sheet to "Include #0"
# - Refresh the Metadata sheet: new Date, new Contact, insert a Jurisdiction row

$wb = $excel.ActiveWorkbook

# --- Rename the second sheet (the "Include" / codes sheet) ---
$codesSheet = $wb.Worksheets.Item("Include from Ferlab.bio CodeS")
$codesSheet.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Date value (row 8, column B)
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# Contact value (row 10, column B)
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# Insert a new row 11 for "Jurisdiction" (pushes Description/Purpose/Copyright/Immutable down by one)
$ws.Rows.Item(11).Insert()

$ws.Range("A11").Value = "Jurisdiction"
# A leading apostrophe forces this to be stored as literal text, giving an explicit
# (empty) string value rather than clearing the cell outright.
$ws.Range("B11").Value = "'"

# Copy the formatting (borders/fill/font/alignment) from the row below (now row 12) onto the
# freshly inserted row 11 so the new row matches the sheet's body style.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
